$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Characters(21, 2).Text = "14"
$ws.Range("C9").Characters(47, 9).Text = "4/7/2024"
$ws.Range("C9").Characters(27, 9).Text = "4/1/2024"

# --- Data cell updates ---
# Cells that change numeric<->text type, fix style via donor PasteSpecial(xlPasteFormats) then set value
$ws.Range("C14").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("F15").Value = "'0"
$excel.CutCopyMode = $false
$ws.Range("C14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("D22").Value = "'0"
$excel.CutCopyMode = $false
$ws.Range("C14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E22").Value = "***.*"
$excel.CutCopyMode = $false
$ws.Range("C16").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Range("D23").Value = 1
$excel.CutCopyMode = $false
$ws.Range("E16").Copy()
$ws.Range("E23").PasteSpecial(-4122)
$ws.Range("E23").Value = -100
$excel.CutCopyMode = $false
$ws.Range("C16").Copy()
$ws.Range("G23").PasteSpecial(-4122)
$ws.Range("G23").Value = 1
$excel.CutCopyMode = $false
$ws.Range("E16").Copy()
$ws.Range("H23").PasteSpecial(-4122)
$ws.Range("H23").Value = 0
$excel.CutCopyMode = $false
$ws.Range("C14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D27").Value = "'0"
$excel.CutCopyMode = $false
$ws.Range("C14").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E27").Value = "***.*"
$excel.CutCopyMode = $false
$ws.Range("C14").Copy()
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("F27").Value = "'0"
$excel.CutCopyMode = $false
$ws.Range("C16").Copy()
$ws.Range("F31").PasteSpecial(-4122)
$ws.Range("F31").Value = 1
$excel.CutCopyMode = $false

# Plain value updates (style/type unchanged)
$ws.Range("H15").Value = -100
$ws.Range("I16").Value = 37
$ws.Range("J16").Value = 47
$ws.Range("K16").Value = -21.27659574468
$ws.Range("L16").Value = -26
$ws.Range("M16").Value = -24.489795918367
$ws.Range("N16").Value = -85.140562248996
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 33.333333333333
$ws.Range("F17").Value = 21
$ws.Range("G17").Value = 23
$ws.Range("H17").Value = -8.695652173913
$ws.Range("I17").Value = 50
$ws.Range("J17").Value = 64
$ws.Range("K17").Value = -21.875
$ws.Range("L17").Value = -10.714285714285
$ws.Range("M17").Value = 78.571428571428
$ws.Range("N17").Value = -15.254237288135
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 9
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = -60
$ws.Range("I18").Value = 55
$ws.Range("J18").Value = 61
$ws.Range("K18").Value = -9.83606557377
$ws.Range("L18").Value = -25.675675675675
$ws.Range("M18").Value = -38.888888888888
$ws.Range("N18").Value = -86.873508353222
$ws.Range("C19").Value = 14
$ws.Range("D19").Value = 18
$ws.Range("E19").Value = -22.222222222222
$ws.Range("F19").Value = 74
$ws.Range("G19").Value = 87
$ws.Range("H19").Value = -14.942528735632
$ws.Range("I19").Value = 249
$ws.Range("J19").Value = 267
$ws.Range("K19").Value = -6.741573033707
$ws.Range("L19").Value = 5.95744680851
$ws.Range("M19").Value = -33.77659574468
$ws.Range("N19").Value = -63.59649122807
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -75
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = -42.857142857142
$ws.Range("I20").Value = 9
$ws.Range("J20").Value = 17
$ws.Range("K20").Value = -47.058823529411
$ws.Range("L20").Value = -25
$ws.Range("M20").Value = 28.571428571428
$ws.Range("N20").Value = -97.345132743362
$ws.Range("C21").Value = 24
$ws.Range("D21").Value = 37
$ws.Range("E21").Value = -35.135135135135
$ws.Range("F21").Value = 112
$ws.Range("G21").Value = 150
$ws.Range("H21").Value = -25.333333333333
$ws.Range("I21").Value = 402
$ws.Range("J21").Value = 458
$ws.Range("K21").Value = -12.227074235807
$ws.Range("L21").Value = -7.586206896551
$ws.Range("M21").Value = -27.173913043478
$ws.Range("N21").Value = -77.107061503416
$ws.Range("F22").Value = 3
$ws.Range("G22").Value = 5
$ws.Range("H22").Value = -40
$ws.Range("I22").Value = 19
$ws.Range("K22").Value = -17.391304347826
$ws.Range("L22").Value = -34.482758620689
$ws.Range("M22").Value = -13.636363636363
$ws.Range("J23").Value = 4
$ws.Range("K23").Value = -50
$ws.Range("D24").Value = 32
$ws.Range("E24").Value = 59.375
$ws.Range("F24").Value = 254
$ws.Range("G24").Value = 140
$ws.Range("H24").Value = 81.428571428571
$ws.Range("I24").Value = 832
$ws.Range("J24").Value = 505
$ws.Range("K24").Value = 64.752475247524
$ws.Range("L24").Value = 39.130434782608
$ws.Range("M24").Value = 80.08658008658
$ws.Range("D25").Value = 21
$ws.Range("E25").Value = 104.761904761905
$ws.Range("F25").Value = 223
$ws.Range("G25").Value = 114
$ws.Range("H25").Value = 95.614035087719
$ws.Range("I25").Value = 710
$ws.Range("J25").Value = 378
$ws.Range("K25").Value = 87.830687830687
$ws.Range("L25").Value = 63.594470046082
$ws.Range("C26").Value = 8
$ws.Range("D26").Value = 12
$ws.Range("E26").Value = -33.333333333333
$ws.Range("F26").Value = 46
$ws.Range("G26").Value = 41
$ws.Range("H26").Value = 12.195121951219
$ws.Range("I26").Value = 148
$ws.Range("J26").Value = 122
$ws.Range("K26").Value = 21.311475409836
$ws.Range("L26").Value = 14.728682170542
$ws.Range("M26").Value = 46.534653465346
$ws.Range("H27").Value = -100
$ws.Range("C28").Value = 4
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = 15
$ws.Range("G28").Value = 11
$ws.Range("H28").Value = 36.363636363636
$ws.Range("I28").Value = 30
$ws.Range("J28").Value = 31
$ws.Range("K28").Value = -3.225806451612
$ws.Range("L28").Value = 25
$ws.Range("I31").Value = 4
$ws.Range("K31").Value = 300
$ws.Range("L31").Value = -50
